$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that currently sits at the end
#    of the "Running the Demos" heading paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Find the "Restore NuGet Packages" bullet and insert a new list item,
#    "Clean & Rebuild the solution", right after it. InsertParagraphAfter
#    on the bullet's Range naturally inherits the ListParagraph style and
#    the numbered-list (numId 1) formatting of the preceding bullet.
$findRange = $d.Content
$findRange.Find.Execute("Restore NuGet Packages", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$restorePara = $findRange.Paragraphs(1)
$restorePara.Range.InsertParagraphAfter()

$newPara = $restorePara.Next()
$newPara.Range.Text = "Clean & Rebuild the solution"

# 3. Re-create the "_GoBack" bookmark, now wrapping the newly inserted
#    paragraph - this mirrors what Word does to mark the last edit
#    location after typing new content.
$bmRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
